$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")
$ws.Activate()

# Update the item's quantity-related figures in row 2
$ws.Range("A2").Value = 2993253
$ws.Range("D2").Value = 221

# Leave the selection where the user last clicked
$ws.Range("B41").Select()
